# The "grandes regiões e unidades da federação" sub-header row (row 6) had
# no data of its own - it was a leftover label row. Correcting the data
# entry removes that empty label row entirely, shifting every row below it
# (7..37) up by one, and dropping the now-unused shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(6).Delete()
